$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row: "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410" ---
# Header row is row 1; the sheet has 21 used columns (A:U) per the dimension/used range.
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count

for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = [string]$cell.Value2

    if ($header.EndsWith("_old")) {
        $base = $header.Substring(0, $header.Length - [string]"_old".Length)
        $cell.Value2 = "$($base)_FV2404"
    }
    elseif ($header.EndsWith("_new")) {
        $base = $header.Substring(0, $header.Length - [string]"_new".Length)
        $cell.Value2 = "$($base)_FV2410"
    }
    # Any other header (e.g. "diff") is left as-is.
}

# --- 2) Turn the data range into an Excel Table (adds autofilter + tableParts) ---
$tableRange = $ws.Range($usedRange.Address(0, 0))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ---
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
